$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'28.183.25"
$ws.Range("E2").Formula = "'  -1.47%  "
$ws.Range("D3").Formula = "'1.805.43"
$ws.Range("E3").Formula = "'  +0.67%  "
$ws.Range("D4").Formula = "'1.001"
$ws.Range("E4").Formula = "'  -0.10%  "
$ws.Range("D5").Formula = "'317.00"
$ws.Range("E5").Formula = "'  +0.97%  "
$ws.Range("D6").Formula = "'1.000"
$ws.Range("E6").Formula = "'  -0.14%  "
$ws.Range("D7").Formula = "'0.5335"
$ws.Range("E7").Formula = "'  -0.42%  "
$ws.Range("E8").Formula = "'  -1.12%  "
$ws.Range("D9").Formula = "'0.07479"
$ws.Range("E9").Formula = "'  -0.71%  "
$ws.Range("D10").Formula = "'42.00"
$ws.Range("E10").Formula = "'  -1.20%  "
$ws.Range("D11").Formula = "'1.097"
$ws.Range("E11").Formula = "'  -2.22%  "
$ws.Range("D12").Formula = "'1.001"
$ws.Range("E12").Formula = "'  -0.08%  "
$ws.Range("D13").Formula = "'6.218"
$ws.Range("E13").Formula = "'  +0.36%  "
$ws.Range("D14").Formula = "'20.54"
$ws.Range("E14").Formula = "'  -2.81%  "
$ws.Range("D15").Formula = "'7.378"
$ws.Range("E15").Formula = "'  -0.59%  "
$ws.Range("D16").Formula = "'1.810.64"
$ws.Range("E16").Formula = "'  +0.95%  "
$ws.Range("D17").Formula = "'89.75"
$ws.Range("E17").Formula = "'  -0.81%  "
$ws.Range("E18").Formula = "'  -0.27%  "
$ws.Range("D19").Formula = "'0.06512"
$ws.Range("E19").Formula = "'  +1.11%  "
$ws.Range("D20").Formula = "'17.43"
$ws.Range("E20").Formula = "'  +0.71%  "
$ws.Range("D21").Formula = "'0.9999"
$ws.Range("E21").Formula = "'  -0.12%  "
$ws.Range("D22").Formula = "'5.927"
$ws.Range("E22").Formula = "'  +0.02%  "
$ws.Range("D23").Formula = "'28.220.26"
$ws.Range("E23").Formula = "'  -1.37%  "
$ws.Range("E24").Formula = "'  -0.27%  "
$ws.Range("D25").Formula = "'2.088"
$ws.Range("E25").Formula = "'  +0.03%  "
$ws.Range("D26").Formula = "'156.35"
$ws.Range("E26").Formula = "'  -3.00%  "
$ws.Range("E27").Formula = "'  -0.43%  "
$ws.Range("D28").Formula = "'2.015.42"
$ws.Range("E28").Formula = "'  +0.77%  "
$ws.Range("D29").Formula = "'2.323"
$ws.Range("E29").Formula = "'  -2.36%  "
$ws.Range("D30").Formula = "'122.11"
$ws.Range("E30").Formula = "'  -0.94%  "
$ws.Range("D31").Formula = "'1.114"
$ws.Range("E31").Formula = "'  -0.90%  "
$ws.Range("D32").Formula = "'0.1095"
$ws.Range("E32").Formula = "'  +7.63%  "
$ws.Range("D33").Formula = "'5.584"
$ws.Range("E33").Formula = "'  -2.37%  "
$ws.Range("D34").Formula = "'3.625"
$ws.Range("E34").Formula = "'  -0.66%  "
$ws.Range("D35").Formula = "'0.07169"
$ws.Range("E35").Formula = "'  +9.04%  "
$ws.Range("D36").Formula = "'0.2226"
$ws.Range("E36").Formula = "'  -3.10%  "
$ws.Range("E37").Formula = "'  -1.15%  "
$ws.Range("D38").Formula = "'5.082"
$ws.Range("E38").Formula = "'  -0.23%  "
$ws.Range("D39").Formula = "'8.450"
$ws.Range("E39").Formula = "'  -2.77%  "
$ws.Range("D40").Formula = "'0.6175"
$ws.Range("E40").Formula = "'  -2.51%  "
$ws.Range("E41").Formula = "'  -3.34%  "
$ws.Range("D42").Formula = "'1.439"
$ws.Range("E42").Formula = "'  +4.20%  "
$ws.Range("D43").Formula = "'1.181"
$ws.Range("E43").Formula = "'  -2.20%  "
$ws.Range("D44").Formula = "'13.36"
$ws.Range("E44").Formula = "'  -1.41%  "
$ws.Range("D45").Formula = "'3.685"
$ws.Range("E45").Formula = "'  +0.38%  "
$ws.Range("D46").Formula = "'0.5768"
$ws.Range("E46").Formula = "'  -2.79%  "
$ws.Range("D47").Formula = "'125.41"
$ws.Range("E47").Formula = "'  -0.44%  "
$ws.Range("D48").Formula = "'1.928"
$ws.Range("E48").Formula = "'  -2.84%  "
$ws.Range("D49").Formula = "'1.185"
$ws.Range("E49").Formula = "'  +1.38%  "
$ws.Range("D50").Formula = "'0.06822"
$ws.Range("E50").Formula = "'  -1.49%  "
$ws.Range("D51").Formula = "'71.84"
$ws.Range("E51").Formula = "'  -1.14%  "
